$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price (D) / Volume(1h) (E) cells hold values such as "43.070.63" or
# "1.00" that Excel would otherwise auto-convert to a Number (or change
# the precision of) when assigned through .Value. Forcing the cell to
# Text format first keeps them stored exactly as the original strings,
# matching the existing inline-string cells on the sheet.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.070.63'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -0.04%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.304.47'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '301.00'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -0.39%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '98.54'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.84%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.522'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +3.76%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.520'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +1.16%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.63'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -1.58%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0792'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.21%  '
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -0.38%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.83'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.38%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.90'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.39%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.663.81'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -0.40%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.331.13'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.09%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.789'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.990.67'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '13.42'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +6.79%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0909'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  +0.59%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.15'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -1.19%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '68.30'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.68%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '239.07'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  +1.22%  '
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -1.59%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.43'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -1.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.79'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.33%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '167.45'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  -0.59%  '
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.44%  '
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -13.39%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '33.37'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.23'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +4.35%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.00'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.03%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +3.19%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '18.10'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +3.74%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -0.57%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -0.39%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -0.78%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.79'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.23%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.98%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.30%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '2.004.00'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -0.33%  '
$ws.Range('B44').Value = 'ApeXProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.13'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -4.35%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.09'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -1.73%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.48'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.70%  '
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -2.98%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '54.51'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.529.82'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -0.47%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.73'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +5.19%  '
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +0.25%  '
